$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Recipes")

# Swap contents of column A ("steps") and column D ("name") across all used rows
# (header row + the 3 data rows), per the diff: A and D columns were transposed.
for ($r = 1; $r -le 4; $r++) {
    $aVal = $ws.Cells.Item($r, 1).Value()
    $dVal = $ws.Cells.Item($r, 4).Value()
    $ws.Cells.Item($r, 1).Value = $dVal
    $ws.Cells.Item($r, 4).Value = $aVal
}
